$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.395.53'
$ws.Range('D3').Value = '3.542.03'
$ws.Range('E3').Value = '  -2.88%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.14'
$ws.Range('E5').Value = '  +0.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.81'
$ws.Range('E6').Value = '  -1.92%  '
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.533.68'
$ws.Range('E7').Value = '  -2.93%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.612'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E10').Value = '  -3.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.82'
$ws.Range('E11').Value = '  -0.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.585'
$ws.Range('E12').Value = '  -3.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '47.45'
$ws.Range('E13').Value = '  -2.33%  '
$ws.Range('E14').Value = '  -3.97%  '
$ws.Range('D15').Value = '4.109.11'
$ws.Range('E15').Value = '  -2.92%  '
$ws.Range('E16').Value = '  -3.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '629.45'
$ws.Range('E17').Value = '  -6.15%  '
$ws.Range('D18').Value = '3.548.15'
$ws.Range('E18').Value = '  -2.56%  '
$ws.Range('D19').Value = '69.418.08'
$ws.Range('E19').Value = '  -1.82%  '
$ws.Range('E20').Value = '  +1.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.41'
$ws.Range('E21').Value = '  -2.10%  '
$ws.Range('E22').Value = '  -1.84%  '
$ws.Range('E23').Value = '  -4.13%  '
$ws.Range('E24').Value = '  -6.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '97.33'
$ws.Range('E25').Value = '  -3.18%  '
$ws.Range('E26').Value = '  -2.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.83'
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('E29').Value = '  -5.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.34'
$ws.Range('E30').Value = '  -6.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.91'
$ws.Range('E31').Value = '  -5.69%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.59'
$ws.Range('E32').Value = '  -4.53%  '
$ws.Range('B33').Value = 'Stacks'
$ws.Range('C33').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.16'
$ws.Range('E33').Value = '  -5.96%  '
$ws.Range('E34').Value = '  -3.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.01'
$ws.Range('E35').Value = '  -3.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '640.10'
$ws.Range('E36').Value = '  +9.31%  '
$ws.Range('E37').Value = '  -2.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.52'
$ws.Range('E38').Value = '  -11.77%  '
$ws.Range('E39').Value = '  -3.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '57.28'
$ws.Range('E40').Value = '  -1.59%  '
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0456'
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('E43').Value = '  -3.60%  '
$ws.Range('D44').Value = '3.390.23'
$ws.Range('E44').Value = '  -5.42%  '
$ws.Range('E45').Value = '  -4.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '32.84'
$ws.Range('E47').Value = '  -5.74%  '
$ws.Range('E48').Value = '  -5.38%  '
$ws.Range('E49').Value = '  -6.66%  '
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '132.57'
$ws.Range('E51').Value = '  -1.86%  '
